$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the duplicate empty paragraph that precedes "The Income Statement"
#    heading (two blank paragraphs had accumulated; only one should remain -
#    the one carrying the bold/underline/size-28 rPr used by the heading
#    paragraphs).
# ---------------------------------------------------------------------------
$dupIdx = -1
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pNext = $d.Paragraphs.Item($i + 1)
    if ($p.Range.Text.Trim() -eq "" -and $pNext.Range.Text.Trim() -eq "") {
        $pNextNext = $d.Paragraphs.Item($i + 2)
        if ($pNextNext.Range.Text -like "The Income Statement*") {
            $dupIdx = $i
            break
        }
    }
}
if ($dupIdx -ge 1) {
    $d.Paragraphs.Item($dupIdx).Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) "Revenue: $60,000/year" (Income Statement section) becomes
#    "Revenue/Cash Flow: $60,000/year"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Revenue: `$", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Revenue/Cash Flow: `$", 2)

# ---------------------------------------------------------------------------
# 3) Under "The Cash Flow Statement" heading, add a new line right before
#    "Cash: $20,000":
#        Revenue/Cash Flow: $60,000/year
#    where "Cash Flow: $60,000/year" is bold.
# ---------------------------------------------------------------------------
$headingIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "The Cash Flow Statement*") {
        $headingIdx = $i
        break
    }
}

$cashIdx = -1
if ($headingIdx -ge 1) {
    $candidate = $d.Paragraphs.Item($headingIdx + 1)
    if ($candidate.Range.Text -like "Cash: `$20,000*") {
        $cashIdx = $headingIdx + 1
    }
}

if ($cashIdx -ge 1) {
    $cashPara = $d.Paragraphs.Item($cashIdx)
    $cashPara.Range.InsertParagraphBefore()

    $newPara = $d.Paragraphs.Item($cashIdx)
    $newPara.Range.Text = "Revenue/Cash Flow: `$60,000/year"

    # Re-fetch the paragraph (text assignment can invalidate cached ranges)
    # and bold the "Cash Flow: $60,000/year" portion only.
    $newPara = $d.Paragraphs.Item($cashIdx)
    $prefixLen = "Revenue/".Length
    $boldStart = $newPara.Range.Start + $prefixLen
    $boldEnd = $newPara.Range.End - 1
    $boldRange = $d.Range($boldStart, $boldEnd)
    $boldRange.Bold = 1
}
